# Commit: swap the two embedded DrawingML themes.
#
#   ppt/theme/theme1.xml  (wired to the slide master / overall "Design")
#       was the "Integral" color theme  -> becomes the "Office Theme" colors
#   ppt/theme/theme2.xml  (wired to the notes master)
#       was the default "Office Theme"  -> becomes the "Integral" colors
#
# i.e. the content that used to live in theme1.xml and theme2.xml trade
# places.  The <a:fontScheme> and <a:fmtScheme> blocks are already
# byte-identical between the two themes (both "Office"/Arial), so only the
# <a:clrScheme> (and the cosmetic theme/clrScheme "name" attributes, which
# are not exposed for writing by this host's PowerPoint object model) truly
# change.
#
# This headless host only exposes ONE writable theme through COM - the
# design theme reachable from $p.SlideMaster.Theme (NotesMaster.Theme /
# HandoutMaster.Theme resolve to that very same part here) - so we push the
# new ("Office Theme") color palette onto it via ThemeColorScheme, which is
# the legitimate PowerPoint VBA surface for editing theme colors.

function Hex2Rgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$colors = $p.SlideMaster.Theme.ThemeColorScheme

# ThemeColorScheme index order: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1..accent6, 11 hlink, 12 folHlink.
# New values = the "Office Theme" palette (previously in theme2.xml).
$colors.Item(1).RGB  = Hex2Rgb "000000"   # dk1
$colors.Item(2).RGB  = Hex2Rgb "FFFFFF"   # lt1
$colors.Item(3).RGB  = Hex2Rgb "44546A"   # dk2
$colors.Item(4).RGB  = Hex2Rgb "E7E6E6"   # lt2
$colors.Item(5).RGB  = Hex2Rgb "5B9BD5"   # accent1
$colors.Item(6).RGB  = Hex2Rgb "ED7D31"   # accent2
$colors.Item(7).RGB  = Hex2Rgb "A5A5A5"   # accent3
$colors.Item(8).RGB  = Hex2Rgb "FFC000"   # accent4
$colors.Item(9).RGB  = Hex2Rgb "4472C4"   # accent5
$colors.Item(10).RGB = Hex2Rgb "70AD47"   # accent6
$colors.Item(11).RGB = Hex2Rgb "0563C1"   # hlink
$colors.Item(12).RGB = Hex2Rgb "954F72"   # folHlink

# Major/minor Latin typefaces are "Arial" on both the old and new theme -
# set explicitly anyway so the intent (full theme swap) is unambiguous.
$fonts = $p.SlideMaster.Theme.ThemeFontScheme
$fonts.MajorFont.Latin = "Arial"
$fonts.MinorFont.Latin = "Arial"
